# io_matrix.xlsx — add a "debug" uart_tx/uart_rx column pair ahead of the
# existing testmode1/testmode2 columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# User selects column F (the rightmost populated column) before inserting
# three new blank columns in front of it (E:G), which pushes the existing
# testmode1 / testmode2 data from E:F out to H:I.
[void]$ws.Range("F1:F1048576").Select()
[void]$ws.Range("E1:G1").EntireColumn.Insert()

# Fill in the newly inserted "debug" signal pair.
$ws.Range("G1").Value = "debug"
$ws.Range("D3").Value = "uart_tx"
$ws.Range("E2").Value = "uart_rx"
$ws.Range("E3").Value = "uart_rx"
$ws.Range("G3").Value = "debug[1]"
